$d = $word.ActiveDocument

# --- 1) Merge the split "Gi" + "t" + "H" + "ub" runs into a single
#        "GitHub" run while keeping the Hyperlink character style. ---
$rng = $d.Content
$rng.Find.Execute("GitHub", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ghStart = $rng.Start
# The replacement text is identical to the existing (multi-run) text, so a
# direct same-text assignment is treated as a no-op; go through a
# temporary placeholder first to force the runs to coalesce.
$rng.Text = "TEMPPLACEHOLDER"
$rng3 = $d.Range($ghStart, $ghStart + 15)
$rng3.Text = "GitHub"
$rng2 = $d.Range($ghStart, $ghStart + 6)
$rng2.Style = "Hyperlink"

# --- 2) Merge the split "Design (" + "fein" + "): Architektur_" + "Fein"
#        + ".png" runs into a single plain run. ---
$rng4 = $d.Content
$rng4.Find.Execute("Design (fein): Architektur_Fein.png", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dfStart = $rng4.Start
$rng4.Text = "TEMPPLACEHOLDERXXXXXXXXXXXXXXXXXXXX"
$rng5 = $d.Range($dfStart, $dfStart + 36)
$rng5.Text = "Design (fein): Architektur_Fein.png"

# --- 3) Add a new list paragraph "Programm: .exe" right after the
#        "Releasenote: Releasenote.docx" entry. ---
$rng6 = $d.Content
$rng6.Find.Execute("Releasenote: Releasenote.docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng6.Collapse(0)
$rng6.InsertParagraphAfter()

$releaseParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Releasenote: Releasenote.docx`r") {
        $releaseParaIndex = $i
        break
    }
}
$newPara = $d.Paragraphs.Item($releaseParaIndex + 1)
$newPara.Range.Text = "Programm: .exe"
